# Update ASV_rank (column H) values on Sheet1 to match the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3   = 960
    11  = 713
    17  = 442
    21  = 795
    31  = 822
    37  = 915
    38  = 670
    47  = 505
    51  = 836
    62  = 859
    70  = 1055
    72  = 860
    80  = 985
    95  = 965
    100 = 863
    112 = 742
    138 = 739
    139 = 829
    140 = 913
    141 = 961
    144 = 1029
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 8).Value = $updates[$row]
}
